$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh adds one new price-report pair ("1a plateado" / "2a
# plateado", Provincia de Melipilla, dated 2021-11-05 / serial 44505) at
# the top of the historical log (row 265), pushing the existing rows
# (265-296) down by two rows (they become 267-298).
$ws.Range("A265:T266").EntireRow.Insert()

# New row 265: "1a plateado"
$ws.Cells.Item(265, 1).Value2 = 11
$ws.Cells.Item(265, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(265, 3).Value = "Bíobío"
$ws.Cells.Item(265, 4).Value2 = 44505
$ws.Cells.Item(265, 5).Value2 = 8
$ws.Cells.Item(265, 6).Value = "Fruta"
$ws.Cells.Item(265, 7).Value2 = 100102
$ws.Cells.Item(265, 8).Value = "Cítricos"
$ws.Cells.Item(265, 9).Value2 = 100102003
$ws.Cells.Item(265, 10).Value = "Limón"
$ws.Cells.Item(265, 11).Value = "Sin especificar"
$ws.Cells.Item(265, 12).Value = "1a plateado"
$ws.Cells.Item(265, 13).Value2 = 600
$ws.Cells.Item(265, 14).Value2 = 7500
$ws.Cells.Item(265, 15).Value2 = 8000
$ws.Cells.Item(265, 16).Value2 = 7750
$ws.Cells.Item(265, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(265, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(265, 19).Value2 = 484
$ws.Cells.Item(265, 20).Value2 = 16

# New row 266: "2a plateado"
$ws.Cells.Item(266, 1).Value2 = 11
$ws.Cells.Item(266, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(266, 3).Value = "Bíobío"
$ws.Cells.Item(266, 4).Value2 = 44505
$ws.Cells.Item(266, 5).Value2 = 8
$ws.Cells.Item(266, 6).Value = "Fruta"
$ws.Cells.Item(266, 7).Value2 = 100102
$ws.Cells.Item(266, 8).Value = "Cítricos"
$ws.Cells.Item(266, 9).Value2 = 100102003
$ws.Cells.Item(266, 10).Value = "Limón"
$ws.Cells.Item(266, 11).Value = "Sin especificar"
$ws.Cells.Item(266, 12).Value = "2a plateado"
$ws.Cells.Item(266, 13).Value2 = 300
$ws.Cells.Item(266, 14).Value2 = 6500
$ws.Cells.Item(266, 15).Value2 = 6500
$ws.Cells.Item(266, 16).Value2 = 6500
$ws.Cells.Item(266, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(266, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(266, 19).Value2 = 406
$ws.Cells.Item(266, 20).Value2 = 16
